# Applies the cryptos-list refresh described in the commit: updates Price (D)
# and Volume(1h) (E) columns for rows 2-51, plus the OKB/Mantle row swap at 45-46.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.279.68"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "3.391.22"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.62"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +9.75%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.391.26"
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("E10").Value = "  +4.45%  "
$ws.Range("E11").Value = "  +4.07%  "
$ws.Range("E12").Value = "  +4.03%  "
$ws.Range("D13").Value = "3.970.42"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("D16").Value = "3.393.10"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.24"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.83%  "
$ws.Range("D18").Value = "61.438.75"
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.07"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +8.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.46"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.80"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "375.98"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("E23").Value = "  +3.45%  "
$ws.Range("D24").Value = "3.524.67"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.77"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000118"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +13.34%  "
$ws.Range("E28").Value = "  +23.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.76"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +13.40%  "
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.15"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.86%  "
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("E33").Value = "  +4.91%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "3.423.00"
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.45"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.59"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +8.99%  "
$ws.Range("E38").Value = "  +7.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.96"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "163.08"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0791"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.55%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.21"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +14.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.43"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.36%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.762"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.43"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.61"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.47"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.99"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +6.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.10"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +15.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.901"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.95%  "
